$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tbl = $sh.Table

$tbl.Cell(1, 1).Shape.TextFrame.TextRange.Text = "T2"
$tbl.Cell(1, 2).Shape.TextFrame.TextRange.Text = "Table 2"
